$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows for ET and SVR classifiers (dt plots for cartpole)
$ws.Range("A7").Value = "ET"
$ws.Range("B7").Value = 40.130000000000003
$ws.Range("C7").Value = 23.62
$ws.Range("D7").Value = 2.36
$ws.Range("E7").Value = 1210.93
$ws.Range("F7").Value = 1937.74
$ws.Range("G7").Value = 193.77
$ws.Range("H7").Value = 1

$ws.Range("A8").Value = "SVR"
$ws.Range("B8").Value = 29.48
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 0.7
$ws.Range("E8").Value = 682
$ws.Range("F8").Value = 182.32
$ws.Range("G8").Value = 18.23
$ws.Range("H8").Value = 0
